$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Row 98: 3/6/2020 ---
$row98 = $tbl.ListRows.Add()
$row98.Range.Cells.Item(1,1).Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$row98.Range.Cells.Item(1,2).Value = 3
$row98.Range.Cells.Item(1,3).Value = 6
$row98.Range.Cells.Item(1,4).Value = 2020
$row98.Range.Cells.Item(1,5).Value = 2630
$row98.Range.Cells.Item(1,6).Value = 1221

# --- Row 99: 4/6/2020 ---
$row99 = $tbl.ListRows.Add()
$row99.Range.Cells.Item(1,1).Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$row99.Range.Cells.Item(1,2).Value = 4
$row99.Range.Cells.Item(1,3).Value = 6
$row99.Range.Cells.Item(1,4).Value = 2020
$row99.Range.Cells.Item(1,5).Value = 2864
$row99.Range.Cells.Item(1,6).Value = 1208

# Apply the same formatting used by the existing table rows (column A uses
# style index 6 - centered, light fill; columns B:D use style index 5 -
# centered) by copying from the row immediately above the newly added ones.
$ws.Range("A97").Copy()
$ws.Range("A98:A99").PasteSpecial(-4122)

$ws.Range("B97:D97").Copy()
$ws.Range("B98:D99").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Match the saved selection/active cell state from the authored workbook.
$ws.Range("G99").Select()
